$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency market data: coin name, link, price, volume(1h)
# Row 9 (OKB) is dropped, all subsequent rows shift up by one, and a new
# row 51 (Cronos) is appended at the bottom.
$rows = @(
    @{ Row = 2; B = 'Bitcoin'; C = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'; D = '30.226.08'; E = '  +0.20%  ' }
    @{ Row = 3; B = 'Ethereum'; C = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'; D = '1.871.57'; E = '  +1.24%  ' }
    @{ Row = 4; B = 'TetherUSD'; C = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'; D = '1.000'; E = '  -0.03%  ' }
    @{ Row = 5; B = 'BNB'; C = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'; D = '234.83'; E = '  -0.13%  ' }
    @{ Row = 6; B = 'USDC'; C = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'; D = '1.000'; E = '  -0.02%  ' }
    @{ Row = 7; B = 'XRP'; C = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'; D = '0.4703'; E = '  +0.22%  ' }
    @{ Row = 8; B = 'Cardano'; C = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; D = '0.2848'; E = '  -1.46%  ' }
    @{ Row = 9; B = 'Dogecoin'; C = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; D = '0.06563'; E = '  +0.33%  ' }
    @{ Row = 10; B = 'Solana'; C = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D = '21.25'; E = '  -2.47%  ' }
    @{ Row = 11; B = 'TRON'; C = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D = '0.07801'; E = '  -1.84%  ' }
    @{ Row = 12; B = 'Litecoin'; C = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D = '96.45'; E = '  -1.03%  ' }
    @{ Row = 13; B = 'WrappedEther'; C = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D = '1.854.11'; E = '  +0.23%  ' }
    @{ Row = 14; B = 'Polygon'; C = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; D = '0.6902'; E = '  +2.35%  ' }
    @{ Row = 15; B = 'Polkadot'; C = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D = '5.095'; E = '  +0.20%  ' }
    @{ Row = 16; B = 'BitcoinCash'; C = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D = '266.96'; E = '  +0.24%  ' }
    @{ Row = 17; B = 'WrappedBTC'; C = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D = '30.208.41'; E = '  +0.24%  ' }
    @{ Row = 18; B = 'Avalanche'; C = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D = '13.70'; E = '  +0.66%  ' }
    @{ Row = 19; B = 'ShibaInu'; C = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D = '0.000007737'; E = '  +1.30%  ' }
    @{ Row = 20; B = 'Dai'; C = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'; D = '1.000'; E = '  +0.00%  ' }
    @{ Row = 21; B = 'WrappedliquidstakedEther2.0'; C = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; D = '2.079.92'; E = '  -0.75%  ' }
    @{ Row = 22; B = 'BinanceUSD'; C = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D = '1.000'; E = '  -0.07%  ' }
    @{ Row = 23; B = 'Uniswap'; C = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D = '5.248'; E = '  +0.96%  ' }
    @{ Row = 24; B = 'Chainlink'; C = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D = '6.161'; E = '  +0.51%  ' }
    @{ Row = 25; B = 'Cosmos'; C = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D = '9.511'; E = '  +4.04%  ' }
    @{ Row = 26; B = 'Monero'; C = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D = '165.70'; E = '  -0.36%  ' }
    @{ Row = 27; B = 'EthereumClassic'; C = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D = '18.75'; E = '  -0.17%  ' }
    @{ Row = 28; B = 'LidoDAOToken'; C = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D = '1.936'; E = '  +0.55%  ' }
    @{ Row = 29; B = 'Toncoin'; C = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D = '1.373'; E = '  -0.68%  ' }
    @{ Row = 30; B = 'Stellar'; C = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D = '0.09926'; E = '  +1.05%  ' }
    @{ Row = 31; B = 'Filecoin'; C = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D = '4.354'; E = '  +1.93%  ' }
    @{ Row = 32; B = 'PancakeSwap'; C = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; D = '1.456'; E = '  -0.53%  ' }
    @{ Row = 33; B = 'InternetComputer(DFINITY)'; C = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D = '4.050'; E = '  +1.44%  ' }
    @{ Row = 34; B = 'Hedera'; C = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; D = '0.04740'; E = '  +1.14%  ' }
    @{ Row = 35; B = 'ARBITRUM'; C = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D = '1.130'; E = '  +1.34%  ' }
    @{ Row = 36; B = 'ImmutableX'; C = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D = '0.7004'; E = '  +0.47%  ' }
    @{ Row = 37; B = 'HuobiToken'; C = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; D = '2.716'; E = '  +0.39%  ' }
    @{ Row = 38; B = 'VeChain'; C = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D = '0.01863'; E = '  -0.09%  ' }
    @{ Row = 39; B = 'MXToken'; C = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D = '2.777'; E = '  +7.06%  ' }
    @{ Row = 40; B = 'FraxShare'; C = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D = '6.244'; E = '  -1.35%  ' }
    @{ Row = 41; B = 'Aave'; C = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'; D = '72.66'; E = '  -0.75%  ' }
    @{ Row = 42; B = 'RenderToken'; C = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D = '1.942'; E = '  +0.74%  ' }
    @{ Row = 43; B = 'PaxDollar'; C = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; D = '1.000'; E = '  +0.10%  ' }
    @{ Row = 44; B = 'TheSandbox'; C = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D = '0.4153'; E = '  +0.81%  ' }
    @{ Row = 45; B = 'TrustWalletToken'; C = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D = '0.8324'; E = '  -0.54%  ' }
    @{ Row = 46; B = 'Quant'; C = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; D = '102.90'; E = '  -0.17%  ' }
    @{ Row = 47; B = 'Maker'; C = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'; D = '981.30'; E = '  +4.58%  ' }
    @{ Row = 48; B = 'Aptos'; C = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D = '7.083'; E = '  +1.44%  ' }
    @{ Row = 49; B = 'EnergySwap'; C = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D = '9.155'; E = '  +0.35%  ' }
    @{ Row = 50; B = 'Elrond'; C = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'; D = '34.50'; E = '  +2.15%  ' }
    @{ Row = 51; B = 'Cronos'; C = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; D = '0.05644'; E = '  -0.22%  ' }
)

foreach ($item in $rows) {
    $r = $item.Row
    $rng = $ws.Range("B$r`:E$r")
    # Force text formatting first so numeric-looking strings (e.g. "1.000",
    # "0.06563") are stored as text, matching the source data, not coerced
    # into numbers/dates by Excel's automatic type detection.
    $rng.NumberFormat = "@"
    $ws.Range("B$r").Value = $item.B
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D
    $ws.Range("E$r").Value = $item.E
    # Drop the temporary text-format override so the cell style matches
    # the original (unstyled) data cells.
    $rng.ClearFormats()
}
